$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "1.00", "0.0000182",
# "67.189.29"); the source data stores them as literal text, so force the
# Text number format before assigning to avoid Excel auto-converting them
# to actual numbers (which would drop formatting like trailing zeros or
# switch to scientific notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.189.29"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.508.63"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.26"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.50"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("E9").Value = "  +7.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.117.38"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.40"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.150.96"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.504.28"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.97"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.70"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000123"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.29"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.09"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.39"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.62"
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.44"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.896"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.76"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.54"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.03"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.811.09"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.82"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0312"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.71"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.49"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.52"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.852"
$ws.Range("E51").Value = "  +0.15%  "
